$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11-14 down to 12-15
$ws.Rows.Item(11).Insert()

# Fill in the new row 11 with the new measurement data
$ws.Cells.Item(11, 1).Value = 1
$ws.Cells.Item(11, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(11, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(11, 4).Value = 44574
$ws.Cells.Item(11, 5).Value = 15
$ws.Cells.Item(11, 6).Value = "Fruta"
$ws.Cells.Item(11, 7).Value = 100103
$ws.Cells.Item(11, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(11, 9).Value = 100103002
$ws.Cells.Item(11, 10).Value = "Ciruela"
$ws.Cells.Item(11, 11).Value = "Black Amber"
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 300
$ws.Cells.Item(11, 14).Value = 18000
$ws.Cells.Item(11, 15).Value = 19000
$ws.Cells.Item(11, 16).Value = 18500
$ws.Cells.Item(11, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(11, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(11, 19).Value = 1028
$ws.Cells.Item(11, 20).Value = 18
